$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Move the "Data das cotações" / "Horário do relatório" label+value
# block from rows 20-21 up to rows 6-7. We use Range.Copy(destination)
# so that the distinctive cell formatting (styles) used by C20/D20
# travels along with the move, instead of re-creating new styles.
# ------------------------------------------------------------------
$ws.Range("C20").Copy($ws.Range("C6"))
$ws.Range("D20").Copy($ws.Range("D6"))
$ws.Range("C21").Copy($ws.Range("C7"))
$ws.Range("D21").Copy($ws.Range("D7"))

# The old location loses its special formatting once the block is
# moved away - reset C20/D20 back to the plain style used by the rest
# of that (otherwise empty) row, by copying the formatting of a
# neighboring plain cell onto them.
$ws.Range("A20").Copy($ws.Range("C20"))
$ws.Range("A20").Copy($ws.Range("D20"))

# Now that formatting has been taken care of, drop the leftover values
# from the old location.
$ws.Range("C20:D21").ClearContents()

# Update the report date / time shown under the new position. The date
# text looks like a date to Excel's auto-conversion, so a direct
# Value assignment would turn it into a date serial number; instead we
# write it as a text formula and then paste-special as a value so the
# cell keeps storing plain text (and keeps its existing style/format).
$ws.Range("C7").Formula = '="12/02/2022"'
$ws.Range("C7").Copy()
$ws.Range("C7").PasteSpecial(-4163)
$ws.Range("D7").Value = "16:48"

# ------------------------------------------------------------------
# Refresh the currency table: only three currencies remain (Dólar,
# Dólar Australiano, Dólar Canadense) with new quoted values; every
# other currency row is wiped out.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Dólar"
$ws.Range("B2").Value = 5.25
$ws.Range("C2").Value = "USD"
$ws.Range("D2").Value = "$"

$ws.Range("A3").Value = "Dólar Australiano"
$ws.Range("B3").Value = 3.75
$ws.Range("C3").Value = "AUD"
$ws.Range("D3").Value = "$"

$ws.Range("A4").Value = "Dólar Canadense"
$ws.Range("B4").Value = 4.13
$ws.Range("C4").Value = "CAD"
$ws.Range("D4").Value = "$"

# Row 5 held "Dólar" (old), now cleared entirely.
$ws.Range("A5:D5").ClearContents()

# Rows 6 and 7, columns A and B (currency name / value) are cleared;
# their C/D columns were already populated above with the moved block.
$ws.Range("A6:B6").ClearContents()
$ws.Range("A7:B7").ClearContents()

# Remaining old currency rows (Euro down through Yuan) are all cleared.
$ws.Range("A8:D18").ClearContents()

Write-Host "Edit complete"
